# fix: Reduce zoom level of excel file
#
# The sheet was zoomed in to 120%; bring it back down to the normal 100%
# zoom level, which also lets row 4 (whose height had been stretched to fit
# the larger on-screen text) shrink back down to the same compact height
# used by the other single-line rows on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the on-screen zoom of the active sheet from 120% back to 100%.
$excel.ActiveWindow.Zoom = 100

# Row 4 no longer needs the extra height that the 120% zoom required;
# match it to the same height already used by the other plain rows (e.g. row 5).
$ws.Rows.Item(4).RowHeight = 14.25
